$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 114, shifting rows 114:193 down to 115:194
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the new data
$ws.Cells.Item(114, 1).Value = 11
$ws.Cells.Item(114, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(114, 3).Value = "Bíobío"
$ws.Cells.Item(114, 4).Value = 44978
$ws.Cells.Item(114, 5).Value = 8
$ws.Cells.Item(114, 6).Value = 100112032
$ws.Cells.Item(114, 7).Value = "Zapallo italiano"
$ws.Cells.Item(114, 8).Value = "Sin especificar"
$ws.Cells.Item(114, 9).Value = "Primera"
$ws.Cells.Item(114, 10).Value = 100
$ws.Cells.Item(114, 11).Value = 7000
$ws.Cells.Item(114, 12).Value = 7500
$ws.Cells.Item(114, 13).Value = 7250
$ws.Cells.Item(114, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(114, 15).Value = "Región Metropolitana"
$ws.Cells.Item(114, 16).Value = 145
$ws.Cells.Item(114, 17).Value = 50
$ws.Cells.Item(114, 18).Value = "Hortaliza"
